$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "ProjectGroup" column to Table1 ---
$lo = $ws.ListObjects.Item(1)
$col = $lo.ListColumns.Add()

# Name the new column via its header cell (same effect as typing a name
# into the header cell of an Excel Table, which is how table columns get
# their display name)
$col.Range.Cells.Item(1, 1).Value = "ProjectGroup"

# --- Fill in the data for the new column ---
$dataRange = $col.DataBodyRange
$dataRange.Cells.Item(1, 1).Value = 1
$dataRange.Cells.Item(2, 1).Value = 2
$dataRange.Cells.Item(3, 1).Value = 3
$dataRange.Cells.Item(4, 1).Value = 4

# --- Clear out the stray ProjectDependency value that no longer applies ---
$ws.Range("C3").ClearContents()

# --- Update the selected cell, matching the saved selection state ---
[void]$ws.Range("E7").Select()
